$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Export Worksheet")
$ws2 = $wb.Worksheets.Item("SQL")

# Update header G1
$ws1.Range("G1").Value = "End PM"

# Update G2:G5 to formatted inline strings (and drop the numeric style
# that used to be applied to these cells, matching the target workbook)
$ws1.Range("G2").Value = "End PM:    4.100"
$ws1.Range("G3").Value = "End PM:    1.600"
$ws1.Range("G4").Value = "End PM:    6.330"
$ws1.Range("G5").Value = "End PM:    7.077"
$ws1.Range("G2:G5").ClearFormats()

# Update SQL text in sheet2 A2
$sql = "select a.ea, a.treatment, a.county, a.route, a.year, ('Beg PM: ' || to_char(a.beg_pm, 990.999)) as `"Beg PM`", ('End PM: ' || to_char(a.end_pm, 990.999)) as `"End PM`", (a.end_pm-a.beg_pm) as length, a.budget_group from s1383currentl a `nwhere a.county = 'SF' `nunion  `nselect b.ea, b.treatment, b.county, b.route, b.year, ('Beg PM: ' || to_char(b.beg_pm, 990.999)) as `"Beg PM`",  ('End PM: ' || to_char(b.end_pm, 990.999)) as `"End PM`",  (b.end_pm-b.beg_pm) as length, b.budget_group from s1383historyl b `nwhere b.county = 'SF' `norder by year"
$ws2.Range("A2").Value = $sql
